$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# --- Row 1: drop the stray repeated "value" header cells in C1:F1 ---
$ws.Range("C1:F1").ClearContents()

# --- Insert a new row for the "L_curve" parameter right after the
#     (soon to be renamed) "production_function" row ---
$ws.Rows.Item(9).Insert()

# Row 8 label changes from "Model" to "production_function"
# (its value, "Sigmoid", is untouched)
$ws.Range("A8").Value = "production_function"

# New row 9: "L_curve" = 1, formatted like the other numeric
# optimization-parameter rows above (e.g. B2, a scientific-format cell)
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# --- Remove the old "Deletion" row, which (after the insert above)
#     now lives at row 17, right before "simulation_timepoints" ---
[void]$ws.Rows.Item(17).Select()
$ws.Rows.Item(17).Delete()
